$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = '[u''[3, 5]'']'
$ws.Cells.Item(2,2).Value = '[u''[1.0, 7.0]'']'
$ws.Cells.Item(2,3).Value = '[ u''Scale, 2'']'
$ws.Cells.Item(2,4).Value = '[''U'']'
$ws.Cells.Item(2,5).Value = '[''2'']'
$ws.Cells.Item(2,6).Value = '[''one'']'
$ws.Cells.Item(2,7).Value = '[''two'']'
$ws.Cells.Item(2,8).Value = '[u"[''C'', 4]"]'
$ws.Cells.Item(2,9).Value = '[u''different'']'
$ws.Cells.Item(2,10).Value = 1
# Row 3
$ws.Cells.Item(3,1).Value = '[u''[11, 8]'']'
$ws.Cells.Item(3,2).Value = '[u''[11, 8]'']'
$ws.Cells.Item(3,3).Value = '['''']'
$ws.Cells.Item(3,4).Value = '[''D'']'
$ws.Cells.Item(3,5).Value = '[''0'']'
$ws.Cells.Item(3,6).Value = '[''none'']'
$ws.Cells.Item(3,7).Value = '[''none'']'
$ws.Cells.Item(3,8).Value = '[u"[''A'', 4]"]'
$ws.Cells.Item(3,9).Value = '[u''same'']'
$ws.Cells.Item(3,10).Value = 2
# Row 4
$ws.Cells.Item(4,1).Value = '[u''[3.0, 1.0, 6]'']'
$ws.Cells.Item(4,2).Value = '[u''[1, 3, 6]'']'
$ws.Cells.Item(4,3).Value = '[u''Scale, 2'']'
$ws.Cells.Item(4,4).Value = '[''changes'']'
$ws.Cells.Item(4,5).Value = '[''2'']'
$ws.Cells.Item(4,6).Value = '[''three'']'
$ws.Cells.Item(4,7).Value = '[''four'']'
$ws.Cells.Item(4,8).Value = '[u"[''B'', 3]"]'
$ws.Cells.Item(4,9).Value = '[u''different'']'
$ws.Cells.Item(4,10).Value = 3
# Row 5
$ws.Cells.Item(5,1).Value = '[u''[3.0, 4.0, 10]'']'
$ws.Cells.Item(5,2).Value = '[u''[5, 6, 10]'']'
$ws.Cells.Item(5,3).Value = '[u''Different Scale, 2'']'
$ws.Cells.Item(5,4).Value = '[''UU'']'
$ws.Cells.Item(5,5).Value = '[''2'']'
$ws.Cells.Item(5,6).Value = '[''one'']'
$ws.Cells.Item(5,7).Value = '[''two'']'
$ws.Cells.Item(5,8).Value = '[u"[''B'', 3]"]'
$ws.Cells.Item(5,9).Value = '[u''different'']'
$ws.Cells.Item(5,10).Value = 4
# Row 6
$ws.Cells.Item(6,1).Value = '[u''[12, 10, 6]'']'
$ws.Cells.Item(6,2).Value = '[u''[12, 10, 6]'']'
$ws.Cells.Item(6,3).Value = '['''']'
$ws.Cells.Item(6,4).Value = '[''DD'']'
$ws.Cells.Item(6,5).Value = '[''0'']'
$ws.Cells.Item(6,6).Value = '[''none'']'
$ws.Cells.Item(6,7).Value = '[''none'']'
$ws.Cells.Item(6,8).Value = '[u"[''A'', 4]"]'
$ws.Cells.Item(6,9).Value = '[u''same'']'
$ws.Cells.Item(6,10).Value = 5
# Row 7
$ws.Cells.Item(7,1).Value = '[u''[11, 8, 5, 4, 1]'']'
$ws.Cells.Item(7,2).Value = '[u''[13.0, 8, 3.0, 4, 1]'']'
$ws.Cells.Item(7,3).Value = '[u''No Scale, 2'']'
$ws.Cells.Item(7,4).Value = '[''changes'']'
$ws.Cells.Item(7,5).Value = '[''2'']'
$ws.Cells.Item(7,6).Value = '[''one'']'
$ws.Cells.Item(7,7).Value = '[''three'']'
$ws.Cells.Item(7,8).Value = '[ u"[''G'', 3]"]'
$ws.Cells.Item(7,9).Value = '[u''different'']'
$ws.Cells.Item(7,10).Value = 6

[void]$ws.Range("I8").Select()
